# Scheduled runner update: refresh market-price derived columns (H-N) across
# the Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) with newly
# polled values. Column layout per sheet:
#   H currentAveragePrice, I currentAveragePriceNQ, J currentAveragePriceHQ,
#   K LevePriceNQ, L LevePriceHQ, M LeveProfitNQ, N LeveProfitHQ
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 32
$ws.Range("I8").Value = 32.7
$ws.Range("J8").Value = 28.5
$ws.Range("K8").Value = 98.10000000000001
$ws.Range("L8").Value = 85.5
$ws.Range("M8").Value = 40.89999999999999
$ws.Range("N8").Value = -363.5
$ws.Range("H13").Value = 1499
$ws.Range("J13").Value = 1499
$ws.Range("L13").Value = 1499
$ws.Range("N13").Value = -1837
$ws.Range("H17").Value = 1723.3334
$ws.Range("J17").Value = 1723.3334
$ws.Range("L17").Value = 5170.0002
$ws.Range("N17").Value = -5506.0002
$ws.Range("H38").Value = 2551.75
$ws.Range("I38").Value = 2551.75
$ws.Range("K38").Value = 7655.25
$ws.Range("M38").Value = -7283.25
$ws.Range("H39").Value = 272.2
$ws.Range("I39").Value = 194.16667
$ws.Range("J39").Value = 389.25
$ws.Range("K39").Value = 582.50001
$ws.Range("L39").Value = 1167.75
$ws.Range("M39").Value = -286.50001
$ws.Range("N39").Value = -1759.75
$ws.Range("H87").Value = 39999.332
$ws.Range("J87").Value = 39999.332
$ws.Range("L87").Value = 39999.332
$ws.Range("N87").Value = -42495.332
$ws.Range("H90").Value = 39999.332
$ws.Range("J90").Value = 39999.332
$ws.Range("L90").Value = 119997.996
$ws.Range("N90").Value = -132477.996
$ws.Range("H103").Value = 800.5714
$ws.Range("I103").Value = 280
$ws.Range("K103").Value = 840
$ws.Range("M103").Value = -254
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()
$ws.Range("H138").Value = 1664.5
$ws.Range("I138").Value = 996.75
$ws.Range("J138").Value = 3000
$ws.Range("K138").Value = 2990.25
$ws.Range("L138").Value = 9000
$ws.Range("M138").Value = 2149.75
$ws.Range("N138").Value = -19280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 15003
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 15003
$ws.Range("K3").Value = 0
$ws.Range("L3").ClearContents()
$ws.Range("M3").Value = 15003
$ws.Range("N3").Value = -15233
$ws.Range("H32").Value = 11200.842
$ws.Range("I32").Value = 6524.385
$ws.Range("J32").Value = 21333.166
$ws.Range("K32").Value = 6524.385
$ws.Range("L32").Value = 21333.166
$ws.Range("M32").Value = -6237.385
$ws.Range("N32").Value = -21907.166
$ws.Range("H61").Value = 12384
$ws.Range("I61").Value = 8426.538
$ws.Range("K61").Value = 8426.538
$ws.Range("M61").Value = -8214.538
$ws.Range("H97").Value = 391
$ws.Range("I97").Value = 334.44446
$ws.Range("J97").Value = 900
$ws.Range("K97").Value = 334.44446
$ws.Range("L97").Value = 900
$ws.Range("M97").Value = 161.55554
$ws.Range("N97").Value = -1892
$ws.Range("H110").Value = 851.125
$ws.Range("I110").Value = 866.3333
$ws.Range("K110").Value = 866.3333
$ws.Range("M110").Value = 1178.6667
$ws.Range("H132").Value = 3781.6667
$ws.Range("I132").Value = 2138
$ws.Range("J132").Value = 12000
$ws.Range("K132").Value = 6414
$ws.Range("L132").Value = 36000
$ws.Range("M132").Value = -3884
$ws.Range("N132").Value = -41060
$ws.Range("H136").Value = 12384
$ws.Range("I136").Value = 8426.538
$ws.Range("K136").Value = 25279.614
$ws.Range("M136").Value = -22729.614

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 3000
$ws.Range("I33").Value = 3000
$ws.Range("K33").Value = 3000
$ws.Range("M33").Value = -2664
$ws.Range("H119").Value = 84500
$ws.Range("I119").Value = 84500
$ws.Range("K119").Value = 84500
$ws.Range("M119").Value = -79662

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 32.833332
$ws.Range("I7").Value = 10
$ws.Range("J7").Value = 64.8
$ws.Range("K7").Value = 10
$ws.Range("L7").Value = 64.8
$ws.Range("M7").Value = 103
$ws.Range("N7").Value = -290.8
$ws.Range("H31").Value = 3059.818
$ws.Range("I31").Value = 1685
$ws.Range("K31").Value = 1685
$ws.Range("M31").Value = -1390
$ws.Range("H34").Value = 3059.818
$ws.Range("I34").Value = 1685
$ws.Range("K34").Value = 1685
$ws.Range("M34").Value = -1483
$ws.Range("H132").Value = 4493.4287
$ws.Range("I132").Value = 2013.75
$ws.Range("J132").Value = 7799.6665
$ws.Range("K132").Value = 6041.25
$ws.Range("L132").Value = 23398.9995
$ws.Range("M132").Value = -3511.25
$ws.Range("N132").Value = -28458.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").ClearContents()
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = 0
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").ClearContents()
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = 0
$ws.Range("H97").Value = 370
$ws.Range("I97").Value = 370
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1110
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -614
$ws.Range("H132").Value = 1266.3334
$ws.Range("J132").Value = 1999
$ws.Range("L132").Value = 17991
$ws.Range("N132").Value = -23051
$ws.Range("H140").Value = 633.125
$ws.Range("I140").Value = 633.125
$ws.Range("K140").Value = 1899.375
$ws.Range("M140").Value = 3280.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 31.125
$ws.Range("I2").Value = 25.166666
$ws.Range("J2").Value = 49
$ws.Range("K2").Value = 25.166666
$ws.Range("L2").Value = 49
$ws.Range("M2").Value = 87.83333400000001
$ws.Range("N2").Value = -275
$ws.Range("H10").Value = 2379.3333
$ws.Range("I10").Value = 2379.3333
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 2379.3333
$ws.Range("L10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = -2210.3333
$ws.Range("H43").Value = 4012.75
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H70").Value = 7749.5
$ws.Range("J70").Value = 8499
$ws.Range("L70").Value = 8499
$ws.Range("N70").Value = -9039.5
$ws.Range("H73").Value = 7749.5
$ws.Range("J73").Value = 8499
$ws.Range("L73").Value = 8499
$ws.Range("N73").Value = -10371
$ws.Range("H80").Value = 33942.855
$ws.Range("I80").Value = 27425
$ws.Range("J80").Value = 42633.332
$ws.Range("K80").Value = 27425
$ws.Range("L80").Value = 42633.332
$ws.Range("M80").Value = -26427
$ws.Range("N80").Value = -44629.332
$ws.Range("H83").Value = 33942.855
$ws.Range("I83").Value = 27425
$ws.Range("J83").Value = 42633.332
$ws.Range("K83").Value = 137125
$ws.Range("L83").Value = 213166.66
$ws.Range("M83").Value = -132133
$ws.Range("N83").Value = -223150.66
$ws.Range("H113").Value = 1556.25
$ws.Range("J113").Value = 1750
$ws.Range("L113").Value = 1750
$ws.Range("N113").Value = -6090
$ws.Range("H132").Value = 2120.1
$ws.Range("I132").Value = 1525.5
$ws.Range("K132").Value = 4576.5
$ws.Range("M132").Value = -2046.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 91.40909000000001
$ws.Range("I2").Value = 91.40909000000001
$ws.Range("K2").Value = 91.40909000000001
$ws.Range("M2").Value = 20.59090999999999
$ws.Range("H7").Value = 1250
$ws.Range("I7").Value = 1250
$ws.Range("K7").Value = 1250
$ws.Range("M7").Value = -1138
$ws.Range("H40").Value = 898
$ws.Range("I40").Value = 898
$ws.Range("K40").Value = 898
$ws.Range("M40").Value = -762
$ws.Range("H126").Value = 1250
$ws.Range("I126").Value = 1250
$ws.Range("K126").Value = 3750
$ws.Range("M126").Value = -1280
$ws.Range("H132").Value = 5498.6665
$ws.Range("I132").Value = 4914.6665
$ws.Range("K132").Value = 14743.9995
$ws.Range("M132").Value = -12213.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 3401.2856
$ws.Range("I14").Value = 3401.2856
$ws.Range("K14").Value = 3401.2856
$ws.Range("M14").Value = -3233.2856
$ws.Range("H44").Value = 46041
$ws.Range("J44").Value = 46041
$ws.Range("L44").Value = 46041
$ws.Range("N44").Value = -47149
$ws.Range("H58").Value = 28847.166
$ws.Range("I58").Value = 24197.8
$ws.Range("K58").Value = 24197.8
$ws.Range("M58").Value = -23889.8
$ws.Range("H132").Value = 1520.2
$ws.Range("I132").Value = 1486.1428
$ws.Range("J132").Value = 1599.6666
$ws.Range("K132").Value = 4458.428400000001
$ws.Range("L132").Value = 4798.9998
$ws.Range("M132").Value = -1928.428400000001
$ws.Range("N132").Value = -9858.9998

